$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LEILA's balance (row 2, column C: Saldo)
$ws.Range("C2").Value = 58752.93

# Update ANDRE's balance (row 4, column C: Saldo)
$ws.Range("C4").Value = 7000

# Remove the ANA (account 005009922) row entirely; rows below shift up
$ws.Range("A6:C6").EntireRow.Delete()
